# diary update up to 3/1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore classic gridlines display (matches the authoring app's saved view).
$excel.ActiveWindow.DisplayGridlines = $true

# --- Row 23: only the "Reflection" cell's style changes (drop the stray italic
#     explanatory-text variant so it matches the plain wrapped-text style used
#     by the rest of the row). Re-use F23's style (same visual style, s=10).
$ws.Range("F23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# --- Row 24: was a placeholder row -> becomes the 2/21 diary entry.
#     Seed formats from row 23 (date cell + wrapped text cells + mood cell),
#     then overwrite the actual values.
$ws.Range("A23:G23").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)

$ws.Range("A24").Value = 43882
$ws.Range("B24").Value = "2PM-7PM"
$ws.Range("C24").Value = "Team"
$ws.Range("D24").Value = "Describe project architecture, social context, and interesting issues/pull requests"
$ws.Range("E24").Value = "All goals"
$ws.Range("F24").Value = "Our project uses elements of but does not strictly implement a lot of different architectural styles and patterns, making an accurate architectural diagram difficult. On the other hand, the social context for the project is well-documented because of its forum archives and github metrics."
$ws.Range("G24").Value = "Good, finished faster than expected"
$ws.Rows.Item(24).RowHeight = 114.9

# --- Row 25: was a placeholder row -> becomes the 2/27 diary entry.
$ws.Range("A23:G23").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)

$ws.Range("A25").Value = 43888
$ws.Range("B25").Value = "5PM-8PM"
$ws.Range("C25").Value = "Class"
$ws.Range("D25").Value = "Attend lecture"
$ws.Range("E25").Value = "Learned about design patterns"
$ws.Range("F25").Value = "There are a ton of design patterns out there, and it would be pretty useful to learn more of them. Good planning saves a lot of time."
$ws.Range("G25").Value = "Positive"
$ws.Rows.Item(25).RowHeight = 58.5

# --- Row 26: was blank -> becomes the 3/1 diary entry.
$ws.Range("A23:G23").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)

$ws.Range("A26").Value = 43891
$ws.Range("B26").Value = "2:00PM-7:30PM"
$ws.Range("C26").Value = "Team"
$ws.Range("D26").Value = "Complete part of homework"
$ws.Range("E26").Value = "All of homework"
$ws.Range("F26").Value = "<what insight(s) did you gain?>"
$ws.Range("G26").Value = "<how did you feel during the activity?>"
$ws.Rows.Item(26).RowHeight = 15.75

# --- Row 27: was blank -> becomes a fresh blank template row (so the diary
#     keeps one open template row at the bottom, as before).
$ws.Range("A23:G23").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)

$ws.Range("A27").Value = "<what day?>"
$ws.Range("B27").Value = "<what time?>"
$ws.Range("C27").Value = "<as applicable, with whom?>"
$ws.Range("D27").Value = "<what did you want to accomplish?>"
$ws.Range("E27").Value = "<what did you actually accomplish?>"
$ws.Range("F27").Value = "<what insight(s) did you gain?>"
$ws.Range("G27").Value = "<how did you feel during the activity?>"
$ws.Rows.Item(27).RowHeight = 15.75

# Row 27's A cell is a template placeholder (plain text), not a date, so
# drop the date-number-format styling it inherited from the copy above.
$ws.Range("B27").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "<what day?>"

# Move the selection/active cell to match where the author ended up editing.
$ws.Range("B30").Select()
